$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Localized for quest display:
# The 3rd entity ("Scavenge for coins" quest row) no longer needs a prerequisite,
# so isHiddenQuest (D3) goes from true(1) to false(0) and prerequisiteKey (E3) is cleared.
$ws.Range("D3").Value = 0
$ws.Range("E3").Clear()

# Update the selected cell shown in the sheet view.
$ws.Range("E5").Select()
